$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 21 de Marzo de 2020 a las 14:46"

# Update province/city case data (table rows 4-61)
$ws.Cells.Item(4, 1).Value = "Madrid"
$ws.Cells.Item(4, 2).Value = 8921
$ws.Cells.Item(4, 3).Value = 1186
$ws.Cells.Item(4, 4).Value = 6931
$ws.Cells.Item(4, 5).Value = 804
$ws.Cells.Item(5, 1).Value = "Cataluña"
$ws.Cells.Item(5, 2).Value = 4203
$ws.Cells.Item(5, 3).Value = 3
$ws.Cells.Item(5, 4).Value = 4078
$ws.Cells.Item(5, 5).Value = 122
$ws.Cells.Item(6, 1).Value = "Araba/Alava"
$ws.Cells.Item(6, 2).Value = 703
$ws.Cells.Item(6, 3).Value = 21
$ws.Cells.Item(6, 4).Value = 655
$ws.Cells.Item(6, 5).Value = 48
$ws.Cells.Item(7, 1).Value = "Navarra"
$ws.Cells.Item(7, 2).Value = 664
$ws.Cells.Item(7, 3).Value = 2
$ws.Cells.Item(7, 4).Value = 652
$ws.Cells.Item(7, 5).Value = 10
$ws.Cells.Item(8, 1).Value = "Valencia/Valencia"
$ws.Cells.Item(8, 2).Value = 627
$ws.Cells.Item(8, 3).Value = 12
$ws.Cells.Item(8, 4).Value = 600
$ws.Cells.Item(8, 5).Value = 15
$ws.Cells.Item(9, 1).Value = "La Rioja"
$ws.Cells.Item(9, 2).Value = 564
$ws.Cells.Item(9, 3).Value = 13
$ws.Cells.Item(9, 4).Value = 536
$ws.Cells.Item(9, 5).Value = 15
$ws.Cells.Item(10, 1).Value = "Bizkaia/Vizcaya"
$ws.Cells.Item(10, 2).Value = 539
$ws.Cells.Item(10, 3).Value = 21
$ws.Cells.Item(10, 4).Value = 522
$ws.Cells.Item(10, 5).Value = 17
$ws.Cells.Item(11, 1).Value = "Asturias"
$ws.Cells.Item(11, 2).Value = 486
$ws.Cells.Item(11, 3).Value = 12
$ws.Cells.Item(11, 4).Value = 467
$ws.Cells.Item(11, 5).Value = 7
$ws.Cells.Item(12, 1).Value = "Malaga"
$ws.Cells.Item(12, 2).Value = 484
$ws.Cells.Item(12, 3).Value = 72
$ws.Cells.Item(12, 4).Value = 464
$ws.Cells.Item(12, 5).Value = 20
$ws.Cells.Item(13, 1).Value = "Ciudad Real"
$ws.Cells.Item(13, 2).Value = 400
$ws.Cells.Item(13, 3).Value = 8
$ws.Cells.Item(13, 4).Value = 364
$ws.Cells.Item(13, 5).Value = 28
$ws.Cells.Item(14, 1).Value = "Alacant/Alicante"
$ws.Cells.Item(14, 2).Value = 372
$ws.Cells.Item(14, 3).Value = 7
$ws.Cells.Item(14, 4).Value = 348
$ws.Cells.Item(14, 5).Value = 17
$ws.Cells.Item(15, 1).Value = "Toledo"
$ws.Cells.Item(15, 2).Value = 370
$ws.Cells.Item(15, 3).Value = 15
$ws.Cells.Item(15, 4).Value = 336
$ws.Cells.Item(15, 5).Value = 19
$ws.Cells.Item(16, 1).Value = "A Coruña"
$ws.Cells.Item(16, 2).Value = 329
$ws.Cells.Item(16, 3).Value = 5
$ws.Cells.Item(16, 4).Value = 326
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(17, 1).Value = "Albacete"
$ws.Cells.Item(17, 2).Value = 327
$ws.Cells.Item(17, 3).Value = 8
$ws.Cells.Item(17, 4).Value = 291
$ws.Cells.Item(17, 5).Value = 28
$ws.Cells.Item(18, 1).Value = "Granada"
$ws.Cells.Item(18, 2).Value = 289
$ws.Cells.Item(18, 3).Value = 72
$ws.Cells.Item(18, 4).Value = 276
$ws.Cells.Item(18, 5).Value = 13
$ws.Cells.Item(19, 1).Value = "Zaragoza"
$ws.Cells.Item(19, 2).Value = 278
$ws.Cells.Item(19, 3).Value = 0
$ws.Cells.Item(19, 4).Value = 264
$ws.Cells.Item(19, 5).Value = 14
$ws.Cells.Item(20, 1).Value = "Burgos"
$ws.Cells.Item(20, 2).Value = 269
$ws.Cells.Item(20, 3).Value = 27
$ws.Cells.Item(20, 4).Value = 175
$ws.Cells.Item(20, 5).Value = 14
$ws.Cells.Item(21, 1).Value = "Salamanca"
$ws.Cells.Item(21, 2).Value = 265
$ws.Cells.Item(21, 3).Value = 13
$ws.Cells.Item(21, 4).Value = 180
$ws.Cells.Item(21, 5).Value = 15
$ws.Cells.Item(22, 1).Value = "Pontevedra"
$ws.Cells.Item(22, 2).Value = 264
$ws.Cells.Item(22, 3).Value = 5
$ws.Cells.Item(22, 4).Value = 262
$ws.Cells.Item(22, 5).Value = 2
$ws.Cells.Item(23, 1).Value = "Illes Balears"
$ws.Cells.Item(23, 2).Value = 246
$ws.Cells.Item(23, 3).Value = 10
$ws.Cells.Item(23, 4).Value = 232
$ws.Cells.Item(23, 5).Value = 4
$ws.Cells.Item(24, 1).Value = "Sevilla"
$ws.Cells.Item(24, 2).Value = 245
$ws.Cells.Item(24, 3).Value = 72
$ws.Cells.Item(24, 4).Value = 243
$ws.Cells.Item(24, 5).Value = 2
$ws.Cells.Item(25, 1).Value = "Caceres"
$ws.Cells.Item(25, 2).Value = 243
$ws.Cells.Item(25, 3).Value = 2
$ws.Cells.Item(25, 4).Value = 231
$ws.Cells.Item(25, 5).Value = 10
$ws.Cells.Item(26, 1).Value = "Valladolid"
$ws.Cells.Item(26, 2).Value = 241
$ws.Cells.Item(26, 3).Value = 13
$ws.Cells.Item(26, 4).Value = 193
$ws.Cells.Item(26, 5).Value = 10
$ws.Cells.Item(27, 1).Value = "Guadalajara"
$ws.Cells.Item(27, 2).Value = 237
$ws.Cells.Item(27, 3).Value = 2
$ws.Cells.Item(27, 4).Value = 231
$ws.Cells.Item(27, 5).Value = 4
$ws.Cells.Item(28, 1).Value = "Gipuzkoa/Guipuzcoa"
$ws.Cells.Item(28, 2).Value = 223
$ws.Cells.Item(28, 3).Value = 21
$ws.Cells.Item(28, 4).Value = 217
$ws.Cells.Item(28, 5).Value = 6
$ws.Cells.Item(29, 1).Value = "Murcia"
$ws.Cells.Item(29, 2).Value = 215
$ws.Cells.Item(29, 3).Value = 1
$ws.Cells.Item(29, 4).Value = 213
$ws.Cells.Item(29, 5).Value = 1
$ws.Cells.Item(30, 1).Value = "Cantabria"
$ws.Cells.Item(30, 2).Value = 215
$ws.Cells.Item(30, 3).Value = 11
$ws.Cells.Item(30, 4).Value = 200
$ws.Cells.Item(30, 5).Value = 4
$ws.Cells.Item(31, 1).Value = "Leon"
$ws.Cells.Item(31, 2).Value = 201
$ws.Cells.Item(31, 3).Value = 3
$ws.Cells.Item(31, 4).Value = 156
$ws.Cells.Item(31, 5).Value = 7
$ws.Cells.Item(32, 1).Value = "Tenerife"
$ws.Cells.Item(32, 2).Value = 192
$ws.Cells.Item(32, 3).Value = 4
$ws.Cells.Item(32, 4).Value = 184
$ws.Cells.Item(32, 5).Value = 4
$ws.Cells.Item(33, 1).Value = "Aragon"
$ws.Cells.Item(33, 2).Value = 174
$ws.Cells.Item(33, 3).Value = 0
$ws.Cells.Item(33, 4).Value = 163
$ws.Cells.Item(33, 5).Value = 11
$ws.Cells.Item(34, 1).Value = "Illes Balears*"
$ws.Cells.Item(34, 2).Value = 169
$ws.Cells.Item(34, 3).Value = 6
$ws.Cells.Item(34, 4).Value = 161
$ws.Cells.Item(34, 5).Value = 2
$ws.Cells.Item(35, 1).Value = "Segovia"
$ws.Cells.Item(35, 2).Value = 157
$ws.Cells.Item(35, 3).Value = 9
$ws.Cells.Item(35, 4).Value = 140
$ws.Cells.Item(35, 5).Value = 13
$ws.Cells.Item(36, 1).Value = "Cordoba"
$ws.Cells.Item(36, 2).Value = 143
$ws.Cells.Item(36, 3).Value = 72
$ws.Cells.Item(36, 4).Value = 140
$ws.Cells.Item(36, 5).Value = 3
$ws.Cells.Item(37, 1).Value = "Jaen"
$ws.Cells.Item(37, 2).Value = 138
$ws.Cells.Item(37, 3).Value = 72
$ws.Cells.Item(37, 4).Value = 136
$ws.Cells.Item(37, 5).Value = 2
$ws.Cells.Item(38, 1).Value = "Cadiz"
$ws.Cells.Item(38, 2).Value = 126
$ws.Cells.Item(38, 3).Value = 72
$ws.Cells.Item(38, 4).Value = 126
$ws.Cells.Item(38, 5).Value = 0
$ws.Cells.Item(39, 1).Value = "Soria"
$ws.Cells.Item(39, 2).Value = 119
$ws.Cells.Item(39, 3).Value = 5
$ws.Cells.Item(39, 4).Value = 71
$ws.Cells.Item(39, 5).Value = 6
$ws.Cells.Item(40, 1).Value = "Avila"
$ws.Cells.Item(40, 2).Value = 114
$ws.Cells.Item(40, 3).Value = 14
$ws.Cells.Item(40, 4).Value = 55
$ws.Cells.Item(40, 5).Value = 7
$ws.Cells.Item(41, 1).Value = "Badajoz"
$ws.Cells.Item(41, 2).Value = 111
$ws.Cells.Item(41, 3).Value = 5
$ws.Cells.Item(41, 4).Value = 104
$ws.Cells.Item(41, 5).Value = 2
$ws.Cells.Item(42, 1).Value = "Castello/Castellon"
$ws.Cells.Item(42, 2).Value = 104
$ws.Cells.Item(42, 3).Value = 1
$ws.Cells.Item(42, 4).Value = 102
$ws.Cells.Item(42, 5).Value = 1
$ws.Cells.Item(43, 1).Value = "Cuenca"
$ws.Cells.Item(43, 2).Value = 94
$ws.Cells.Item(43, 3).Value = 5
$ws.Cells.Item(43, 4).Value = 84
$ws.Cells.Item(43, 5).Value = 5
$ws.Cells.Item(44, 1).Value = "Ourense"
$ws.Cells.Item(44, 2).Value = 74
$ws.Cells.Item(44, 3).Value = 5
$ws.Cells.Item(44, 4).Value = 74
$ws.Cells.Item(44, 5).Value = 0
$ws.Cells.Item(45, 1).Value = "Gran Canaria"
$ws.Cells.Item(45, 2).Value = 70
$ws.Cells.Item(45, 3).Value = 0
$ws.Cells.Item(45, 4).Value = 69
$ws.Cells.Item(45, 5).Value = 1
$ws.Cells.Item(46, 1).Value = "Zamora"
$ws.Cells.Item(46, 2).Value = 59
$ws.Cells.Item(46, 3).Value = 3
$ws.Cells.Item(46, 4).Value = 42
$ws.Cells.Item(46, 5).Value = 2
$ws.Cells.Item(47, 1).Value = "Igualada, Vilanova del Cami, Santa Margarida de Montbui y Odena"
$ws.Cells.Item(47, 2).Value = 58
$ws.Cells.Item(47, 3).Value = 0
$ws.Cells.Item(47, 4).Value = 58
$ws.Cells.Item(47, 5).Value = 3
$ws.Cells.Item(48, 1).Value = "Lugo"
$ws.Cells.Item(48, 2).Value = 54
$ws.Cells.Item(48, 3).Value = 5
$ws.Cells.Item(48, 4).Value = 53
$ws.Cells.Item(48, 5).Value = 1
$ws.Cells.Item(49, 1).Value = "Almeria"
$ws.Cells.Item(49, 2).Value = 53
$ws.Cells.Item(49, 3).Value = 72
$ws.Cells.Item(49, 4).Value = 53
$ws.Cells.Item(49, 5).Value = 0
$ws.Cells.Item(50, 1).Value = "Palencia"
$ws.Cells.Item(50, 2).Value = 41
$ws.Cells.Item(50, 3).Value = 2
$ws.Cells.Item(50, 4).Value = 25
$ws.Cells.Item(50, 5).Value = 0
$ws.Cells.Item(51, 1).Value = "Teruel"
$ws.Cells.Item(51, 2).Value = 40
$ws.Cells.Item(51, 3).Value = 0
$ws.Cells.Item(51, 4).Value = 38
$ws.Cells.Item(51, 5).Value = 2
$ws.Cells.Item(52, 1).Value = "Huelva"
$ws.Cells.Item(52, 2).Value = 37
$ws.Cells.Item(52, 3).Value = 72
$ws.Cells.Item(52, 4).Value = 37
$ws.Cells.Item(52, 5).Value = 0
$ws.Cells.Item(53, 1).Value = "Huesca"
$ws.Cells.Item(53, 2).Value = 34
$ws.Cells.Item(53, 3).Value = 0
$ws.Cells.Item(53, 4).Value = 34
$ws.Cells.Item(53, 5).Value = 0
$ws.Cells.Item(54, 1).Value = "Melilla"
$ws.Cells.Item(54, 2).Value = 25
$ws.Cells.Item(54, 3).Value = 0
$ws.Cells.Item(54, 4).Value = 25
$ws.Cells.Item(54, 5).Value = 0
$ws.Cells.Item(55, 1).Value = "Fuerteventura"
$ws.Cells.Item(55, 2).Value = 12
$ws.Cells.Item(55, 3).Value = 0
$ws.Cells.Item(55, 4).Value = 12
$ws.Cells.Item(55, 5).Value = 0
$ws.Cells.Item(56, 1).Value = "La Palma"
$ws.Cells.Item(56, 2).Value = 7
$ws.Cells.Item(56, 3).Value = 0
$ws.Cells.Item(56, 4).Value = 7
$ws.Cells.Item(56, 5).Value = 0
$ws.Cells.Item(57, 1).Value = "Arroyo de la Luz"
$ws.Cells.Item(57, 2).Value = 7
$ws.Cells.Item(57, 3).Value = 0
$ws.Cells.Item(57, 4).Value = 7
$ws.Cells.Item(57, 5).Value = 0
$ws.Cells.Item(58, 1).Value = "Ceuta"
$ws.Cells.Item(58, 2).Value = 5
$ws.Cells.Item(58, 3).Value = 0
$ws.Cells.Item(58, 4).Value = 5
$ws.Cells.Item(58, 5).Value = 0
$ws.Cells.Item(59, 1).Value = "Lanzarote"
$ws.Cells.Item(59, 2).Value = 3
$ws.Cells.Item(59, 3).Value = 0
$ws.Cells.Item(59, 4).Value = 3
$ws.Cells.Item(59, 5).Value = 0
$ws.Cells.Item(60, 1).Value = "La Gomera"
$ws.Cells.Item(60, 2).Value = 3
$ws.Cells.Item(60, 3).Value = 2
$ws.Cells.Item(60, 4).Value = 1
$ws.Cells.Item(60, 5).Value = 0
$ws.Cells.Item(61, 1).Value = "El Hierro"
$ws.Cells.Item(61, 2).Value = 1
$ws.Cells.Item(61, 3).Value = 0
$ws.Cells.Item(61, 4).Value = 1
$ws.Cells.Item(61, 5).Value = 0
